$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 458.5
$ws.Cells.Item(2, 9).Value = 167
$ws.Cells.Item(2, 11).Value = 167
$ws.Cells.Item(2, 13).Value = -54

$ws.Cells.Item(9, 8).Value = 309
$ws.Cells.Item(9, 9).Value = 291.42856
$ws.Cells.Item(9, 11).Value = 291.42856
$ws.Cells.Item(9, 13).Value = -122.42856

$ws.Cells.Item(29, 8).Value = 2026
$ws.Cells.Item(29, 9).Value = 2368
$ws.Cells.Item(29, 10).Value = 1000
$ws.Cells.Item(29, 11).Value = 7104
$ws.Cells.Item(29, 12).Value = 3000
$ws.Cells.Item(29, 13).Value = -6823
$ws.Cells.Item(29, 14).Value = -3562

$ws.Cells.Item(38, 8).Value = 53.625
$ws.Cells.Item(38, 9).Value = 39.857143
$ws.Cells.Item(38, 10).Value = 150
$ws.Cells.Item(38, 11).Value = 119.571429
$ws.Cells.Item(38, 12).Value = 450
$ws.Cells.Item(38, 13).Value = 252.428571
$ws.Cells.Item(38, 14).Value = -1194

$ws.Cells.Item(43, 8).Value = 1865.5555
$ws.Cells.Item(43, 9).Value = 1888.5
$ws.Cells.Item(43, 10).Value = 1682
$ws.Cells.Item(43, 11).Value = 1888.5
$ws.Cells.Item(43, 12).Value = 1682
$ws.Cells.Item(43, 13).Value = -1819.5
$ws.Cells.Item(43, 14).Value = -1820

$ws.Cells.Item(58, 8).Value = 84.25
$ws.Cells.Item(58, 9).Value = 84.25
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 252.75
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -102.75
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 2072.4849
$ws.Cells.Item(132, 9).Value = 2074.75
$ws.Cells.Item(132, 11).Value = 6224.25
$ws.Cells.Item(132, 13).Value = -3694.25

$ws.Cells.Item(137, 8).Value = 5252.2
$ws.Cells.Item(137, 9).Value = 4089.4443
$ws.Cells.Item(137, 10).Value = 6996.3335
$ws.Cells.Item(137, 11).Value = 12268.3329
$ws.Cells.Item(137, 12).Value = 20989.0005
$ws.Cells.Item(137, 13).Value = -9718.332900000001
$ws.Cells.Item(137, 14).Value = -26089.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4431.079
$ws.Cells.Item(32, 9).Value = 2661.1292
$ws.Cells.Item(32, 11).Value = 2661.1292
$ws.Cells.Item(32, 13).Value = -2374.1292

$ws.Cells.Item(61, 8).Value = 2389.4666
$ws.Cells.Item(61, 9).Value = 2294.4614
$ws.Cells.Item(61, 11).Value = 2294.4614
$ws.Cells.Item(61, 13).Value = -2082.4614

$ws.Cells.Item(63, 8).Value = 3089.9092
$ws.Cells.Item(63, 9).Value = 3089.9092
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 3089.9092
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = -2403.9092
$ws.Cells.Item(63, 14).ClearContents()

$ws.Cells.Item(66, 8).Value = 3089.9092
$ws.Cells.Item(66, 9).Value = 3089.9092
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 15449.546
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = -12017.546
$ws.Cells.Item(66, 14).ClearContents()

$ws.Cells.Item(97, 8).Value = 15620
$ws.Cells.Item(97, 10).Value = 26875
$ws.Cells.Item(97, 12).Value = 26875
$ws.Cells.Item(97, 14).Value = -27867

$ws.Cells.Item(102, 8).Value = 3802.244
$ws.Cells.Item(102, 9).Value = 3690.5938
$ws.Cells.Item(102, 11).Value = 3690.5938
$ws.Cells.Item(102, 13).Value = -2068.5938

$ws.Cells.Item(110, 8).Value = 1121.3
$ws.Cells.Item(110, 9).Value = 1144.1765
$ws.Cells.Item(110, 10).Value = 991.6667
$ws.Cells.Item(110, 11).Value = 1144.1765
$ws.Cells.Item(110, 12).Value = 991.6667
$ws.Cells.Item(110, 13).Value = 900.8235
$ws.Cells.Item(110, 14).Value = -5081.6667

$ws.Cells.Item(132, 9).Value = 1760.44
$ws.Cells.Item(132, 10).Value = 4248.5
$ws.Cells.Item(132, 11).Value = 5281.32
$ws.Cells.Item(132, 12).Value = 12745.5
$ws.Cells.Item(132, 13).Value = -2751.32
$ws.Cells.Item(132, 14).Value = -17805.5

$ws.Cells.Item(136, 8).Value = 2389.4666
$ws.Cells.Item(136, 9).Value = 2294.4614
$ws.Cells.Item(136, 11).Value = 6883.3842
$ws.Cells.Item(136, 13).Value = -4333.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2394.7317
$ws.Cells.Item(94, 9).Value = 524.0345
$ws.Cells.Item(94, 11).Value = 524.0345
$ws.Cells.Item(94, 13).Value = -73.03449999999998

$ws.Cells.Item(105, 8).Value = 1454.9615
$ws.Cells.Item(105, 9).Value = 1327.5652
$ws.Cells.Item(105, 10).Value = 2431.6667
$ws.Cells.Item(105, 11).Value = 1327.5652
$ws.Cells.Item(105, 12).Value = 2431.6667
$ws.Cells.Item(105, 13).Value = 419.4348
$ws.Cells.Item(105, 14).Value = -5925.6667

$ws.Cells.Item(107, 8).Value = 2690.2
$ws.Cells.Item(107, 9).Value = 2765.4443
$ws.Cells.Item(107, 11).Value = 2765.4443
$ws.Cells.Item(107, 13).Value = -845.4443000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3120.1177
$ws.Cells.Item(31, 9).Value = 2019.3
$ws.Cells.Item(31, 10).Value = 4692.7144
$ws.Cells.Item(31, 11).Value = 2019.3
$ws.Cells.Item(31, 12).Value = 4692.7144
$ws.Cells.Item(31, 13).Value = -1724.3
$ws.Cells.Item(31, 14).Value = -5282.7144

$ws.Cells.Item(34, 8).Value = 3120.1177
$ws.Cells.Item(34, 9).Value = 2019.3
$ws.Cells.Item(34, 10).Value = 4692.7144
$ws.Cells.Item(34, 11).Value = 2019.3
$ws.Cells.Item(34, 12).Value = 4692.7144
$ws.Cells.Item(34, 13).Value = -1817.3
$ws.Cells.Item(34, 14).Value = -5096.7144

$ws.Cells.Item(62, 8).Value = 1900
$ws.Cells.Item(62, 9).Value = 1350
$ws.Cells.Item(62, 11).Value = 1350
$ws.Cells.Item(62, 13).Value = -726

$ws.Cells.Item(65, 8).Value = 1900
$ws.Cells.Item(65, 9).Value = 1350
$ws.Cells.Item(65, 11).Value = 6750
$ws.Cells.Item(65, 13).Value = -3630

$ws.Cells.Item(68, 8).Value = 49999
$ws.Cells.Item(68, 10).Value = 49999
$ws.Cells.Item(68, 12).Value = 49999
$ws.Cells.Item(68, 14).Value = -51497

$ws.Cells.Item(71, 8).Value = 49999
$ws.Cells.Item(71, 10).Value = 49999
$ws.Cells.Item(71, 12).Value = 149997
$ws.Cells.Item(71, 14).Value = -157485

$ws.Cells.Item(131, 8).Value = 43030.21
$ws.Cells.Item(131, 10).Value = 43030.21
$ws.Cells.Item(131, 12).Value = 43030.21
$ws.Cells.Item(131, 14).Value = -53110.21

$ws.Cells.Item(134, 8).Value = 28116.945
$ws.Cells.Item(134, 9).Value = 37050.81
$ws.Cells.Item(134, 11).Value = 111152.43
$ws.Cells.Item(134, 13).Value = -108617.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 420
$ws.Cells.Item(34, 9).Value = 183.33333
$ws.Cells.Item(34, 11).Value = 549.99999
$ws.Cells.Item(34, 13).Value = -465.99999

$ws.Cells.Item(46, 8).Value = 159835.9
$ws.Cells.Item(46, 9).Value = 417360.88
$ws.Cells.Item(46, 10).Value = 1359
$ws.Cells.Item(46, 11).Value = 1252082.64
$ws.Cells.Item(46, 12).Value = 4077
$ws.Cells.Item(46, 13).Value = -1251991.64
$ws.Cells.Item(46, 14).Value = -4259

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 2017
$ws.Cells.Item(43, 9).Value = 2017
$ws.Cells.Item(43, 11).Value = 2017
$ws.Cells.Item(43, 13).Value = -1866

$ws.Cells.Item(80, 8).Value = 2272.7334
$ws.Cells.Item(80, 9).Value = 1998
$ws.Cells.Item(80, 10).Value = 2341.4167
$ws.Cells.Item(80, 11).Value = 1998
$ws.Cells.Item(80, 12).Value = 2341.4167
$ws.Cells.Item(80, 13).Value = -1000
$ws.Cells.Item(80, 14).Value = -4337.4167

$ws.Cells.Item(83, 8).Value = 2272.7334
$ws.Cells.Item(83, 9).Value = 1998
$ws.Cells.Item(83, 10).Value = 2341.4167
$ws.Cells.Item(83, 11).Value = 9990
$ws.Cells.Item(83, 12).Value = 11707.0835
$ws.Cells.Item(83, 13).Value = -4998
$ws.Cells.Item(83, 14).Value = -21691.0835

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3650
$ws.Cells.Item(68, 9).Value = 3200.25
$ws.Cells.Item(68, 11).Value = 3200.25
$ws.Cells.Item(68, 13).Value = -2451.25

$ws.Cells.Item(71, 8).Value = 3650
$ws.Cells.Item(71, 9).Value = 3200.25
$ws.Cells.Item(71, 11).Value = 16001.25
$ws.Cells.Item(71, 13).Value = -12257.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 9250
$ws.Cells.Item(58, 9).Value = 9250
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 9250
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -8942
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 229816.73
$ws.Cells.Item(132, 9).Value = 2593.8948
$ws.Cells.Item(132, 10).Value = 1668894.6
$ws.Cells.Item(132, 11).Value = 7781.6844
$ws.Cells.Item(132, 12).Value = 5006683.800000001
$ws.Cells.Item(132, 13).Value = -5251.6844
$ws.Cells.Item(132, 14).Value = -5011743.800000001
